$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Funciones ya hechas según el Administrador" list (M2:N9) is being
# extended down to also include the remaining "MODULO PRODUCTOS" items
# (rows 10-13), reusing the same look as the already-listed M9:N9 row.
$ws.Range("M9").Copy()
$ws.Range("M10:N10").Merge()
$ws.Range("M10:N10").PasteSpecial(-4122)
$ws.Range("M10").Value = "2. Mostrar productos"

$ws.Range("M9").Copy()
$ws.Range("M11:N11").Merge()
$ws.Range("M11:N11").PasteSpecial(-4122)
$ws.Range("M11").Value = "3. Buscar productos"

$ws.Range("M9").Copy()
$ws.Range("M12:N12").Merge()
$ws.Range("M12:N12").PasteSpecial(-4122)
$ws.Range("M12").Value = "4. Actualizar productos"

$ws.Range("M9").Copy()
$ws.Range("M13:N13").Merge()
$ws.Range("M13:N13").PasteSpecial(-4122)
$ws.Range("M13").Value = "5. Eliminar productos"

$excel.CutCopyMode = 0

# New "Funciones ya hechas según el Vendedor" header, mirroring the plain
# (unstyled) Administrador header in M2 — Vendedor/Usuario functions are
# still pending (see commit message).
$ws.Range("M16").Value = "Funciones ya hechas según el  Vendedor"
